# chore: adapt column header formatting to respective input file names
#
# Renames the comparison-table headers from the generic "_old"/"_new"
# suffixes to the concrete format-version suffixes "_FV2310"/"_FV2404",
# freezes the header row, and wraps the sheet's data range in a proper
# Excel Table (ListObject) named "Table1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the header row (row 1, columns A:U).
#    A:J  -> "<name>_FV2310"   (was "<name>_old")
#    K    -> "diff"            (unchanged)
#    L:U  -> "<name>_FV2404"   (was "<name>_new")
# ---------------------------------------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2310"
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2404"
}

# ---------------------------------------------------------------------
# 2. Freeze the header row.
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3. Turn the used range into an Excel Table ("Table1") with an
#    autofilter and one column per header.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:U69")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
